$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 128, shifting existing rows 128-210 down to 129-211.
$ws.Rows(128).Insert()

# Populate the newly inserted row 128 with the new record's data.
$ws.Range("A128").Value2 = 11
$ws.Range("B128").Value = "Vega Monumental Concepción"
$ws.Range("C128").Value = "Bíobío"
$ws.Range("D128").Value2 = 44824
$ws.Range("E128").Value2 = 8
$ws.Range("F128").Value2 = 100112003
$ws.Range("G128").Value = "Ajo"
$ws.Range("H128").Value = "Chino"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value2 = 400
$ws.Range("K128").Value2 = 17000
$ws.Range("L128").Value2 = 18000
$ws.Range("M128").Value2 = 17500
$ws.Range("N128").Value = "$/caja 10 kilos"
$ws.Range("O128").Value = "China"
$ws.Range("P128").Value2 = 1750
$ws.Range("Q128").Value2 = 10
$ws.Range("R128").Value = "Hortaliza"
